$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Replace "Emre Abale" with "Rob Oudman" in A2
$ws.Range("A2").Value = "Rob Oudman"

# Update selection to A2
$ws.Range("A2").Select()
